# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the 813e05ad-... entry on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-23 22:46:18"
$wsZhCn.Range("H3").Value = "2016-03-23 22:46:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-23 22:46:22"
$wsDeDe.Range("H3").Value = "2016-03-23 22:47:01"
